$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.277420333333333
$ws.Range("H2").Value = 3.832261
$ws.Range("I2").Value = 0.01913942624337554
$ws.Range("J2").Value = 0.01913942624337554
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 20.424575
$ws.Range("N2").Value = 61.273725
$ws.Range("O2").Value = 0.203732656096709
$ws.Range("P2").Value = 0.2037326560967089
$ws.Range("Q2").Value = 26.09076740469167
$ws.Range("R2").Value = 234.816906642225
$ws.Range("S2").Value = 0.003899326144729956
$ws.Range("T2").Value = 0.003899326144729955
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.277420333333333
$ws.Range("H3").Value = 3.832261
$ws.Range("I3").Value = 0.01913942624337554
$ws.Range("J3").Value = 0.01913942624337554
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 57.16769933333334
$ws.Range("N3").Value = 171.503098
$ws.Range("O3").Value = 0.5702408607336045
$ws.Range("P3").Value = 0.5702408607336045
$ws.Range("Q3").Value = 73.02718153828646
$ws.Range("R3").Value = 657.2446338445781
$ws.Range("S3").Value = 0.01091408289496981
$ws.Range("T3").Value = 0.01091408289496981
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.277420333333333
$ws.Range("H4").Value = 3.832261
$ws.Range("I4").Value = 0.01913942624337554
$ws.Range("J4").Value = 0.01913942624337554
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1182046666666667
$ws.Range("N4").Value = 0.354614
$ws.Range("O4").Value = 0.001179077200040937
$ws.Range("P4").Value = 0.001179077200040937
$ws.Range("Q4").Value = 0.1509970446948889
$ws.Range("R4").Value = 1.358973402254
$ws.Range("S4").Value = 0.00002256686110542926
$ws.Range("T4").Value = 0.00002256686110542926
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.277420333333333
$ws.Range("H5").Value = 3.832261
$ws.Range("I5").Value = 0.01913942624337554
$ws.Range("J5").Value = 0.01913942624337554
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 22.54136766666667
$ws.Range("N5").Value = 67.62410300000001
$ws.Range("O5").Value = 0.2248474059696456
$ws.Range("P5").Value = 0.2248474059696456
$ws.Range("Q5").Value = 28.79480139854256
$ws.Range("R5").Value = 259.153212586883
$ws.Range("S5").Value = 0.004303450342570351
$ws.Range("T5").Value = 0.00430345034257035
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 45.44725166666667
$ws.Range("H6").Value = 136.341755
$ws.Range("I6").Value = 0.6809303864519871
$ws.Range("J6").Value = 0.6809303864519872
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 20.424575
$ws.Range("N6").Value = 61.273725
$ws.Range("O6").Value = 0.203732656096709
$ws.Range("P6").Value = 0.2037326560967089
$ws.Range("Q6").Value = 928.2408002097084
$ws.Range("R6").Value = 8354.167201887376
$ws.Range("S6").Value = 0.1387277562488218
$ws.Range("T6").Value = 0.1387277562488218
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 45.44725166666667
$ws.Range("H7").Value = 136.341755
$ws.Range("I7").Value = 0.6809303864519871
$ws.Range("J7").Value = 0.6809303864519872
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 57.16769933333334
$ws.Range("N7").Value = 171.503098
$ws.Range("O7").Value = 0.5702408607336045
$ws.Range("P7").Value = 0.5702408607336045
$ws.Range("Q7").Value = 2598.114818806333
$ws.Range("R7").Value = 23383.03336925699
$ws.Range("S7").Value = 0.3882943296700471
$ws.Range("T7").Value = 0.3882943296700471
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 45.44725166666667
$ws.Range("H8").Value = 136.341755
$ws.Range("I8").Value = 0.6809303864519871
$ws.Range("J8").Value = 0.6809303864519872
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.1182046666666667
$ws.Range("N8").Value = 0.354614
$ws.Range("O8").Value = 0.001179077200040937
$ws.Range("P8").Value = 0.001179077200040937
$ws.Range("Q8").Value = 5.372077234174444
$ws.Range("R8").Value = 48.34869510757
$ws.Range("S8").Value = 0.0008028694934806019
$ws.Range("T8").Value = 0.0008028694934806021
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 45.44725166666667
$ws.Range("H9").Value = 136.341755
$ws.Range("I9").Value = 0.6809303864519871
$ws.Range("J9").Value = 0.6809303864519872
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 22.54136766666667
$ws.Range("N9").Value = 67.62410300000001
$ws.Range("O9").Value = 0.2248474059696456
$ws.Range("P9").Value = 0.2248474059696456
$ws.Range("Q9").Value = 1024.443209257863
$ws.Range("R9").Value = 9219.988883320766
$ws.Range("S9").Value = 0.1531054310396376
$ws.Range("T9").Value = 0.1531054310396376
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4966396666666666
$ws.Range("H10").Value = 1.489919
$ws.Range("I10").Value = 0.007441088905245192
$ws.Range("J10").Value = 0.007441088905245193
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 20.424575
$ws.Range("N10").Value = 61.273725
$ws.Range("O10").Value = 0.203732656096709
$ws.Range("P10").Value = 0.2037326560967089
$ws.Range("Q10").Value = 10.14365411980833
$ws.Range("R10").Value = 91.292887078275
$ws.Range("S10").Value = 0.001515992806917355
$ws.Range("T10").Value = 0.001515992806917355
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.4966396666666666
$ws.Range("H11").Value = 1.489919
$ws.Range("I11").Value = 0.007441088905245192
$ws.Range("J11").Value = 0.007441088905245193
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 57.16769933333334
$ws.Range("N11").Value = 171.503098
$ws.Range("O11").Value = 0.5702408607336045
$ws.Range("P11").Value = 0.5702408607336045
$ws.Range("Q11").Value = 28.39174714100689
$ws.Range("R11").Value = 255.525724269062
$ws.Range("S11").Value = 0.004243212942122292
$ws.Range("T11").Value = 0.004243212942122294
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.4966396666666666
$ws.Range("H12").Value = 1.489919
$ws.Range("I12").Value = 0.007441088905245192
$ws.Range("J12").Value = 0.007441088905245193
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.1182046666666667
$ws.Range("N12").Value = 0.354614
$ws.Range("O12").Value = 0.001179077200040937
$ws.Range("P12").Value = 0.001179077200040937
$ws.Range("Q12").Value = 0.05870512625177778
$ws.Range("R12").Value = 0.528346136266
$ws.Range("S12").Value = 0.00000877361827165218
$ws.Range("T12").Value = 0.000008773618271652182
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.4966396666666666
$ws.Range("H13").Value = 1.489919
$ws.Range("I13").Value = 0.007441088905245192
$ws.Range("J13").Value = 0.007441088905245193
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 22.54136766666667
$ws.Range("N13").Value = 67.62410300000001
$ws.Range("O13").Value = 0.2248474059696456
$ws.Range("P13").Value = 0.2248474059696456
$ws.Range("Q13").Value = 11.19493732418411
$ws.Range("R13").Value = 100.754435917657
$ws.Range("S13").Value = 0.001673109537933892
$ws.Range("T13").Value = 0.001673109537933892
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 19.52156333333333
$ws.Range("H14").Value = 58.56469
$ws.Range("I14").Value = 0.2924890983993922
$ws.Range("J14").Value = 0.2924890983993922
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 20.424575
$ws.Range("N14").Value = 61.273725
$ws.Range("O14").Value = 0.203732656096709
$ws.Range("P14").Value = 0.2037326560967089
$ws.Range("Q14").Value = 398.7196344189167
$ws.Range("R14").Value = 3588.47670977025
$ws.Range("S14").Value = 0.05958958089623984
$ws.Range("T14").Value = 0.05958958089623984
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 19.52156333333333
$ws.Range("H15").Value = 58.56469
$ws.Range("I15").Value = 0.2924890983993922
$ws.Range("J15").Value = 0.2924890983993922
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 57.16769933333334
$ws.Range("N15").Value = 171.503098
$ws.Range("O15").Value = 0.5702408607336045
$ws.Range("P15").Value = 0.5702408607336045
$ws.Range("Q15").Value = 1116.002863156624
$ws.Range("R15").Value = 10044.02576840962
$ws.Range("S15").Value = 0.1667892352264653
$ws.Range("T15").Value = 0.1667892352264654
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 19.52156333333333
$ws.Range("H16").Value = 58.56469
$ws.Range("I16").Value = 0.2924890983993922
$ws.Range("J16").Value = 0.2924890983993922
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.1182046666666667
$ws.Range("N16").Value = 0.354614
$ws.Range("O16").Value = 0.001179077200040937
$ws.Range("P16").Value = 0.001179077200040937
$ws.Range("Q16").Value = 2.307539886628889
$ws.Range("R16").Value = 20.76785897966
$ws.Range("S16").Value = 0.0003448672271832533
$ws.Range("T16").Value = 0.0003448672271832534
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 19.52156333333333
$ws.Range("H17").Value = 58.56469
$ws.Range("I17").Value = 0.2924890983993922
$ws.Range("J17").Value = 0.2924890983993922
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 22.54136766666667
$ws.Range("N17").Value = 67.62410300000001
$ws.Range("O17").Value = 0.2248474059696456
$ws.Range("P17").Value = 0.2248474059696456
$ws.Range("Q17").Value = 440.0427365247856
$ws.Range("R17").Value = 3960.38462872307
$ws.Range("S17").Value = 0.06576541504950377
$ws.Range("T17").Value = 0.06576541504950377
